$wb = $excel.ActiveWorkbook

# Sheet 1: "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value = 7967
$ws1.Range("F5").Value = 940
$ws1.Range("F6").Value = 292
$ws1.Range("F7").Value = 804
$ws1.Range("F8").Value = 613
$ws1.Range("F9").Value = 94
$ws1.Range("F13").Value = 3267
$ws1.Range("F14").Value = 205
$ws1.Range("F15").Value = 99
$ws1.Range("F17").Value = 754
$ws1.Range("F19").Value = 462
$ws1.Range("F21").Value = 257
$ws1.Range("F22").Value = 228
$ws1.Range("F23").Value = 326
$ws1.Range("F26").Value = 115
$ws1.Range("F27").Value = 280
$ws1.Range("F28").Value = 24
$ws1.Range("F31").Value = 506
$ws1.Range("F32").Value = 541
$ws1.Range("F33").Value = 25
$ws1.Range("F34").Value = 35
$ws1.Range("F35").Value = 11

# Sheet 2: "演出" (Performance) - row 6, G6
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("G6").Value = 408

# Sheet 4: "全部类型" (All types)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F5").Value = 7967
$ws4.Range("F7").Value = 940
$ws4.Range("F8").Value = 292
$ws4.Range("F9").Value = 804
$ws4.Range("F10").Value = 613
$ws4.Range("F11").Value = 94
$ws4.Range("F16").Value = 3267
$ws4.Range("F17").Value = 205
$ws4.Range("F18").Value = 99
$ws4.Range("F21").Value = 754
$ws4.Range("F24").Value = 462
$ws4.Range("F26").Value = 257
$ws4.Range("F27").Value = 228
$ws4.Range("F28").Value = 327
$ws4.Range("F31").Value = 115
$ws4.Range("F32").Value = 280
$ws4.Range("F33").Value = 24
$ws4.Range("F36").Value = 506
$ws4.Range("F37").Value = 541
$ws4.Range("F38").Value = 25
$ws4.Range("F39").Value = 35
$ws4.Range("F40").Value = 12
$ws4.Range("G44").Value = 408
